# Retraining the forecast model for Dragosel Laslea
# - Shift every timestamp in column A (rows 2-97) forward by 2 days.
# - Update column B (Actual Production (MW)) with the refreshed values for
#   rows 2-42; the remaining rows keep their existing (zero) production.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wind production values (row -> MW) for the rows whose reading changed.
$newValues = @{
    2 = 876; 3 = 863; 4 = 856; 5 = 850; 6 = 849; 7 = 818; 8 = 783; 9 = 777;
    10 = 780; 11 = 766; 12 = 751; 13 = 740; 14 = 714; 15 = 696; 16 = 684;
    17 = 687; 19 = 685; 20 = 686; 21 = 694; 22 = 707; 23 = 711; 24 = 732;
    25 = 782; 26 = 824; 27 = 813; 28 = 822; 29 = 813; 30 = 816; 31 = 823;
    32 = 831; 33 = 814; 34 = 771; 35 = 750; 36 = 721; 37 = 694; 38 = 707;
    39 = 696; 40 = 699; 41 = 731; 42 = 749
}

for ($row = 2; $row -le 97; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value2 = $cellA.Value2 + 2

    if ($newValues.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
    }
}
